# Apply "feat: add 2022-Q3 data" edit.
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "总计" (summary) sheet: insert a new top data row for 2022-Q3 and push
#    the existing 2022-Q2 / 2022-Q1 rows down by one.
# ---------------------------------------------------------------------------
$wsTotal = $wb.Worksheets.Item(1)

# Duplicate row 3 (2022-Q1, 4, 0.06) down into row 4 - keeps formatting intact
# and already carries the values the new row 4 needs (only the index column
# A needs to change from 1 -> 2).
$wsTotal.Range("A3:D3").Copy($wsTotal.Range("A4:D4"))
$wsTotal.Range("A4").Value = 2

# Row 3 becomes what used to be row 2's quarter (2022-Q2 / 0.05), row index
# (A3) and count (C3) stay the same.
$wsTotal.Range("B3").Value = "2022-Q2"
$wsTotal.Range("D3").Value = 0.05

# Row 2 becomes the new quarter, 2022-Q3 / 0.04.
$wsTotal.Range("B2").Value = "2022-Q3"
$wsTotal.Range("D2").Value = 0.04

# ---------------------------------------------------------------------------
# 2) Add a brand-new "2022-Q3" sheet (holdings detail), positioned right
#    after "总计" and before "2022-Q2". Easiest reliable way to inherit all
#    number formats/styles is to clone the existing "2022-Q2" sheet (same
#    layout) and then overwrite the cell values.
# ---------------------------------------------------------------------------
$wsQ2 = $wb.Worksheets.Item("2022-Q2")
$wsQ2.Copy($null, $wsTotal)
$wsQ3 = $wb.Worksheets.Item(2)
$wsQ3.Name = "2022-Q3"

# Force text formatting on the columns that must stay literal strings
# (leading zeros / fixed decimal display), matching the source data.
$wsQ3.Range("B2:B5").NumberFormat = "@"
$wsQ3.Range("D2:G5").NumberFormat = "@"

$wsQ3.Range("B2").Value = "008707"
$wsQ3.Range("C2").Value = "建信富时100指数（QDII）美元现汇 A"
$wsQ3.Range("D2").Value = "0.48"
$wsQ3.Range("E2").Value = "89.38"
$wsQ3.Range("F2").Value = "3.24"
$wsQ3.Range("G2").Value = "0.0156"
$wsQ3.Range("H2").Value = 10

$wsQ3.Range("B3").Value = "539003"
$wsQ3.Range("C3").Value = "建信富时100指数（QDII）人民币A"
$wsQ3.Range("D3").Value = "0.48"
$wsQ3.Range("E3").Value = "89.38"
$wsQ3.Range("F3").Value = "3.24"
$wsQ3.Range("G3").Value = "0.0156"
$wsQ3.Range("H3").Value = 10

$wsQ3.Range("B4").Value = "008706"
$wsQ3.Range("C4").Value = "建信富时100指数（QDII）人民币 C"
$wsQ3.Range("D4").Value = "0.19"
$wsQ3.Range("E4").Value = "89.38"
$wsQ3.Range("F4").Value = "3.24"
$wsQ3.Range("G4").Value = "0.0062"
$wsQ3.Range("H4").Value = 10

$wsQ3.Range("B5").Value = "008708"
$wsQ3.Range("C5").Value = "建信富时100指数（QDII）美元现汇 C"
$wsQ3.Range("D5").Value = "0.19"
$wsQ3.Range("E5").Value = "89.38"
$wsQ3.Range("F5").Value = "3.24"
$wsQ3.Range("G5").Value = "0.0062"
$wsQ3.Range("H5").Value = 10

# ---------------------------------------------------------------------------
# 3) Restore the originally-selected sheet ("2022-Q1", now the 4th tab) as
#    the active one, since copying a sheet switches focus to the copy.
# ---------------------------------------------------------------------------
$wb.Worksheets.Item("2022-Q1").Select()
